$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.561.78"
$ws.Range("E2").Value = "  -2.57%  "
$ws.Range("D3").Value = "2.003.89"
$ws.Range("E3").Value = "  -4.25%  "
$ws.Range("E4").Value = "  +1.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "329.90"
$ws.Range("E5").Value = "  -3.97%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.013"
$ws.Range("E6").Value = "  +1.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5008"
$ws.Range("E7").Value = "  -4.43%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4219"
$ws.Range("E8").Value = "  -4.71%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "54.57"
$ws.Range("E9").Value = "  -0.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09036"
$ws.Range("E10").Value = "  -3.24%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.118"
$ws.Range("E11").Value = "  -4.37%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "23.32"
$ws.Range("E12").Value = "  -6.28%  "
$ws.Range("D13").Value = "2.038.01"
$ws.Range("E13").Value = "  -2.42%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.046"
$ws.Range("E14").Value = "  -6.44%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.469"
$ws.Range("E15").Value = "  -6.31%  "
$ws.Range("E16").Value = "  +1.21%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "94.50"
$ws.Range("E17").Value = "  -6.73%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001115"
$ws.Range("E18").Value = "  -3.99%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06677"
$ws.Range("E19").Value = "  +0.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.67"
$ws.Range("E20").Value = "  -7.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.013"
$ws.Range("E21").Value = "  +1.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.975"
$ws.Range("E22").Value = "  -5.82%  "
$ws.Range("D23").Value = "29.601.46"
$ws.Range("E23").Value = "  -2.33%  "
$ws.Range("E24").Value = "  -4.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.303"
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.78"
$ws.Range("E26").Value = "  -2.53%  "
$ws.Range("E27").Value = "  -5.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.338"
$ws.Range("E28").Value = "  -6.69%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.296"
$ws.Range("E29").Value = "  -8.88%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "128.27"
$ws.Range("E30").Value = "  -3.61%  "
$ws.Range("E31").Value = "  -7.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09954"
$ws.Range("E32").Value = "  -4.80%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.568"
$ws.Range("E33").Value = "  -6.42%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.836"
$ws.Range("E34").Value = "  -6.61%  "
$ws.Range("E35").Value = "  -1.59%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02468"
$ws.Range("E36").Value = "  -6.35%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "9.274"
$ws.Range("E37").Value = "  -9.45%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06410"
$ws.Range("E38").Value = "  -6.46%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.305"
$ws.Range("E39").Value = "  -3.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6540"
$ws.Range("E40").Value = "  -6.49%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "11.68"
$ws.Range("E41").Value = "  -6.79%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.2047"
$ws.Range("E42").Value = "  -7.47%  "
$ws.Range("E43").Value = "  +1.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6347"
$ws.Range("E44").Value = "  -7.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.44"
$ws.Range("E45").Value = "  -6.33%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.192"
$ws.Range("E46").Value = "  -6.86%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.305"
$ws.Range("E47").Value = "  -5.14%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.508"
$ws.Range("E48").Value = "  -3.64%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000340"
$ws.Range("E49").Value = "  -0.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06989"
$ws.Range("E50").Value = "  -3.37%  "
$ws.Range("E51").Value = "  -6.69%  "
